# Weekly crime data refresh: advance the report one week
# (Volume 29 Number 50, week of 12/12/2022-12/18/2022)
# to (Volume 29 Number 51, week of 12/19/2022-12/25/2022), and update
# the Week to Date / 28 Day / Year to Date / 2 Year crime-complaint
# figures in the "Crime Complaints" table, plus the Historical
# Perspective labels (unaffected values) on rows 36-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/number and reporting week ---
$ws.Range("A8").Value = "Volume 29   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# --- Crime Complaints table (rows 14-30) ---
# Columns: C=WTD 2022, D=WTD 2021, E=WTD %Chg, F=28Day 2022, G=28Day 2021,
#          H=28Day %Chg, I=YTD 2022, J=YTD 2021, K=YTD %Chg, L=2Yr %Chg,
#          M=12Yr %Chg, N=29Yr %Chg

# Row 14: Murder
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "0"
$ws.Range("E14").Value = "***.*"
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 15
$ws.Range("J14").Value = 16
$ws.Range("K14").Value = -6.25
$ws.Range("L14").Value = 25
$ws.Range("M14").Value = -31.818181818181
$ws.Range("N14").Value = -57.142857142857

# Row 15: Rape
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 45
$ws.Range("J15").Value = 41
$ws.Range("K15").Value = 9.756097560975
$ws.Range("L15").Value = 7.142857142857
$ws.Range("M15").Value = 18.421052631578
$ws.Range("N15").Value = -42.307692307692

# Row 16: Robbery
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = 22.222222222222
$ws.Range("F16").Value = 34
$ws.Range("G16").Value = 45
$ws.Range("H16").Value = -24.444444444444
$ws.Range("I16").Value = 454
$ws.Range("J16").Value = 374
$ws.Range("K16").Value = 21.39037433155
$ws.Range("L16").Value = 43.670886075949
$ws.Range("M16").Value = 5.092592592592
$ws.Range("N16").Value = -67.525035765379

# Row 17: Fel. Assault
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = -22.222222222222
$ws.Range("F17").Value = 56
$ws.Range("G17").Value = 61
$ws.Range("H17").Value = -8.196721311475
$ws.Range("I17").Value = 721
$ws.Range("J17").Value = 722
$ws.Range("K17").Value = -0.138504155124
$ws.Range("L17").Value = 6.499261447562
$ws.Range("M17").Value = 68.457943925233
$ws.Range("N17").Value = -14.775413711583

# Row 18: Burglary
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = 13.043478260869
$ws.Range("I18").Value = 296
$ws.Range("J18").Value = 246
$ws.Range("K18").Value = 20.325203252032
$ws.Range("L18").Value = 24.894514767932
$ws.Range("M18").Value = -20.855614973262
$ws.Range("N18").Value = -83.807439824945

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -21.276595744680
$ws.Range("I19").Value = 746
$ws.Range("J19").Value = 528
$ws.Range("K19").Value = 41.287878787878
$ws.Range("L19").Value = 74.707259953161
$ws.Range("M19").Value = 142.207792207792
$ws.Range("N19").Value = 39.179104477611

# Row 20: G.L.A.
$ws.Range("C20").Value = 18
$ws.Range("D20").Value = 20
$ws.Range("E20").Value = -10
$ws.Range("F20").Value = 51
$ws.Range("G20").Value = 62
$ws.Range("H20").Value = -17.741935483871
$ws.Range("I20").Value = 452
$ws.Range("J20").Value = 565
$ws.Range("K20").Value = -20
$ws.Range("L20").Value = 62.007168458781
$ws.Range("M20").Value = 26.610644257703
$ws.Range("N20").Value = -71.661442006269

# Row 21: TOTAL
$ws.Range("C21").Value = 55
$ws.Range("D21").Value = 62
$ws.Range("E21").Value = -11.290322580645
$ws.Range("F21").Value = 209
$ws.Range("G21").Value = 244
$ws.Range("H21").Value = -14.344262295082
$ws.Range("I21").Value = 2729
$ws.Range("J21").Value = 2492
$ws.Range("K21").Value = 9.510433386837
$ws.Range("L21").Value = 37.135678391959
$ws.Range("M21").Value = 39.305768249106
$ws.Range("N21").Value = -56.792273590880

# Row 22: Transit
$ws.Range("D22").Value = 1
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -40
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = 56
$ws.Range("N22").Value = "***.*"

# Row 23: Housing
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -83.333333333333
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = -60
$ws.Range("I23").Value = 107
$ws.Range("J23").Value = 121
$ws.Range("K23").Value = -11.570247933884
$ws.Range("L23").Value = 7
$ws.Range("M23").Value = 44.594594594594
$ws.Range("N23").Value = "***.*"

# Row 24: Petit Larceny
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 34.782608695652
$ws.Range("F24").Value = 142
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = 77.5
$ws.Range("I24").Value = 1489
$ws.Range("J24").Value = 1080
$ws.Range("K24").Value = 37.870370370370
$ws.Range("L24").Value = 52.249488752556
$ws.Range("M24").Value = 97.480106100795
$ws.Range("N24").Value = "***.*"

# Row 25: Misd. Assault
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -17.647058823529
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 88
$ws.Range("H25").Value = -35.227272727272
$ws.Range("I25").Value = 856
$ws.Range("J25").Value = 934
$ws.Range("K25").Value = -8.351177730192
$ws.Range("L25").Value = -14.570858283433
$ws.Range("M25").Value = -10.553814002089
$ws.Range("N25").Value = "***.*"

# Row 26: UCR Rape*
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "***.*"
$ws.Range("F26").Value = 5
$ws.Range("H26").Value = 25
$ws.Range("I26").Value = 76
$ws.Range("K26").Value = 31.034482758620
$ws.Range("L26").Value = 31.034482758620

# Row 27: Other Sex Crimes
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -14.285714285714
$ws.Range("I27").Value = 71
$ws.Range("J27").Value = 82
$ws.Range("K27").Value = -13.414634146341
$ws.Range("L27").Value = 4.411764705882

# Row 28: Shooting Vic.
$ws.Range("C28").Value = 2
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 400
$ws.Range("I28").Value = 47
$ws.Range("K28").Value = -34.722222222222
$ws.Range("L28").Value = -7.843137254901
$ws.Range("M28").Value = -27.692307692307
$ws.Range("N28").Value = -64.122137404580

# Row 29: Shooting Inc.
$ws.Range("C29").Value = 2
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 400
$ws.Range("I29").Value = 44
$ws.Range("K29").Value = -25.423728813559
$ws.Range("L29").Value = -4.347826086956
$ws.Range("M29").Value = -16.981132075471
$ws.Range("N29").Value = -64.227642276422

# Row 30: Hate Crimes - values unchanged, nothing to do.
